$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Re-order the worker identity data (ID doc + name) across rows 16-20.
# New order derived from the authoritative shared-strings table:
#   row16: 1143329013 / HANDER OMAR MANRIQUE ZABALA
#   row17: 80241179   / RAMON NICOLAS NAVARRO BALLESTAS
#   row18: 1001968151 / MANUEL JOSE GUZMAN DE LA ROSA   (unchanged)
#   row19: 73163482   / FRANKLIN MANUEL DIAZ MUÑOZ
#   row20: 33334065   / INGRIS MARINA DIAZ MUÑOZ

$ws.Range("C16").Value = "1143329013"
$ws.Range("D16").Value = "HANDER OMAR MANRIQUE ZABALA"

$ws.Range("C17").Value = "80241179"
$ws.Range("D17").Value = "RAMON NICOLAS NAVARRO BALLESTAS"

$ws.Range("C18").Value = "1001968151"
$ws.Range("D18").Value = "MANUEL JOSE GUZMAN DE LA ROSA"

$ws.Range("C19").Value = "73163482"
$ws.Range("D19").Value = "FRANKLIN MANUEL DIAZ MUÑOZ"

$ws.Range("C20").Value = "33334065"
$ws.Range("D20").Value = "INGRIS MARINA DIAZ MUÑOZ"

# Update the mora value for row 19
$ws.Range("G19").Value = 737717
